# Apply updated crypto price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.172.40"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "1.829.82"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").Value = "'313.31"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").Value = "'0.4692"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").Value = "'0.3664"
$ws.Range("E8").Value = "  -0.49%  "

$ws.Range("D9").Value = "'0.07404"
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").Value = "'0.8805"
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("D11").Value = "'20.34"
$ws.Range("E11").Value = "  -0.14%  "

$ws.Range("D12").Value = "1.883.76"
$ws.Range("E12").Value = "  +2.22%  "

$ws.Range("D13").Value = "'0.07671"
$ws.Range("E13").Value = "  +5.21%  "

$ws.Range("D14").Value = "'5.388"
$ws.Range("E14").Value = "  -1.29%  "

$ws.Range("D15").Value = "'93.16"
$ws.Range("E15").Value = "  +0.66%  "

$ws.Range("D16").Value = "'6.536"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").Value = "'0.000008731"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("D20").Value = "27.585.54"
$ws.Range("E20").Value = "  +1.74%  "

$ws.Range("D21").Value = "'14.63"
$ws.Range("E21").Value = "  -0.89%  "

$ws.Range("D22").Value = "'5.247"
$ws.Range("E22").Value = "  -1.14%  "

$ws.Range("D23").Value = "'10.63"
$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("D24").Value = "2.089.20"
$ws.Range("E24").Value = "  +1.75%  "

$ws.Range("D25").Value = "'1.880"
$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("D26").Value = "'151.42"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").Value = "'2.117"
$ws.Range("E28").Value = "  -1.77%  "

$ws.Range("D29").Value = "'5.188"
$ws.Range("E29").Value = "  -1.49%  "

$ws.Range("D30").Value = "'116.70"
$ws.Range("E30").Value = "  -0.40%  "

$ws.Range("D31").Value = "'0.08931"
$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").Value = "'0.7457"
$ws.Range("E32").Value = "  -1.67%  "

$ws.Range("D33").Value = "'1.164"
$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("E34").Value = "  +1.17%  "

$ws.Range("D35").Value = "'4.515"
$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("D36").Value = "'1.010"
$ws.Range("E36").Value = "  +0.32%  "

$ws.Range("D37").Value = "'2.539"
$ws.Range("E37").Value = "  +4.92%  "

$ws.Range("D38").Value = "'1.093"
$ws.Range("E38").Value = "  -0.86%  "

$ws.Range("D39").Value = "'0.05299"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("D40").Value = "'0.01938"
$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("D41").Value = "'7.307"
$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("D42").Value = "'2.940"
$ws.Range("E42").Value = "  -1.78%  "

$ws.Range("D43").Value = "'0.5276"
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("D44").Value = "'0.1643"
$ws.Range("E44").Value = "  -1.10%  "

$ws.Range("D45").Value = "'8.398"
$ws.Range("E45").Value = "  -1.66%  "

$ws.Range("D46").Value = "'0.4911"
$ws.Range("E46").Value = "  -0.60%  "

$ws.Range("D47").Value = "'10.46"
$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("D48").Value = "'1.010"
$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("D49").Value = "'104.43"
$ws.Range("E49").Value = "  +0.68%  "

$ws.Range("D50").Value = "'1.654"
$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("D51").Value = "'0.06277"
$ws.Range("E51").Value = "  -0.49%  "
